$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the sheet from "Sheet1" to "misc."
$ws.Name = "misc."

# Add new header cells K4/L4 ("other_indexes" / "commodity")
$ws.Range("K4").Value = "other_indexes"
$ws.Range("L4").Value = "commodity"

# New row 11: flo_emis rule for gas CCS processes -> co2captured commodity
$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95
$ws.Range("K11").Value = "co2"
$ws.Range("L11").Value = "co2captured"

# New row 12: flo_emis rule for coal/oil CCS processes -> co2captured commodity
$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85
$ws.Range("K12").Value = "co2"
$ws.Range("L12").Value = "co2captured"

# Column width adjustments (column E widened to fit new content, column K newly sized)
$ws.Columns.Item(5).ColumnWidth = 9.498697916666666
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666

# Move/restore the active selection to D13 (below the newly added rows)
$ws.Range("D13").Select()
